$p = $ppt.ActivePresentation

# --- Update the cached "datetimeFigureOut" date field text wherever it
# appears (slide master's Date Placeholder + every slide layout's Date
# Placeholder), from 12/18/2023 to 12/29/2023. ---
$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "12/29/2023"
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "12/29/2023"
        }
    }
}

# --- Slide 1: resize/rename the "00-build" step box to
# "02-setting-up-a-project", and resize its connecting elbow connector
# to match. ---
$s = $p.Slides.Item(1)

$rect = $s.Shapes.Item("Rectangle 231")
$rect.Left = 288
$rect.Width = 96
$rect.TextFrame.TextRange.Text = "02-setting-up-a-project"

$elbow = $s.Shapes.Item("Elbow Connector 59")
$elbow.Width = 36
